# Generate Report for Archive
# Updates the localization status from "Ready for handoff" to "In Translation"
# on the Overview sheet (zh-cn / de-de status columns) and on the per-locale
# "zh-cn" / "de-de" sheets (Status column), and shrinks the corresponding
# status columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de), row 2 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C), row 2 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C), row 2 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
